# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp text in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 05:05"

# --- Update country stats (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) ---

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 514992   # B - Casos totales
$ws.Cells.Item(5, 3).Value = 143      # C - Nuevos casos
$ws.Cells.Item(5, 5).Value = 279096   # E - Recuperados
$ws.Cells.Item(5, 7).Value = 27       # G - Muertes hoy
$ws.Cells.Item(5, 8).Value = 29341    # H - Muertes

# Row 57: Bolivia
$ws.Cells.Item(57, 2).Value = 9982    # B
$ws.Cells.Item(57, 3).Value = 390     # C
$ws.Cells.Item(57, 4).Value = 968     # D
$ws.Cells.Item(57, 5).Value = 8701    # E
$ws.Cells.Item(57, 7).Value = 3       # G
$ws.Cells.Item(57, 8).Value = 313     # H

# Row 63: Ghana
$ws.Cells.Item(63, 2).Value = 8070    # B
$ws.Cells.Item(63, 3).Value = 189     # C
$ws.Cells.Item(63, 4).Value = 2947    # D
$ws.Cells.Item(63, 5).Value = 5087    # E

# Row 71: Honduras
$ws.Cells.Item(71, 2).Value = 5202    # B
$ws.Cells.Item(71, 3).Value = 108     # C
$ws.Cells.Item(71, 4).Value = 537     # D
$ws.Cells.Item(71, 5).Value = 4453    # E
$ws.Cells.Item(71, 7).Value = 11      # G
$ws.Cells.Item(71, 8).Value = 212     # H

# Row 158: Birmania
$ws.Cells.Item(158, 2).Value = 228    # B
$ws.Cells.Item(158, 3).Value = 4      # C
$ws.Cells.Item(158, 5).Value = 84     # E

# --- Swap order/data of Santa Lucia <-> Belice (rows 201/202) ---
$ws.Cells.Item(201, 1).Value = "Belice"
$ws.Cells.Item(201, 4).Value = 16
$ws.Cells.Item(201, 8).Value = 2

$ws.Cells.Item(202, 1).Value = "Santa Lucia"
$ws.Cells.Item(202, 4).Value = 18
$ws.Cells.Item(202, 8).Value = 0

# --- Swap order/data of Montserrat <-> Seychelles (rows 210/211) ---
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1
